{"js": "// Update the date line and the 25 division problems in the practice\n// worksheet table. Each entry below gives the exact text to search for\n// (matchCase) and the ordered list of replacement strings to apply to\n// each successive match (handles the one duplicated \"36\u00f75=\" cell, whose\n// two occurrences map to two different replacements).\nconst edits = [\n  { find: \"2024-05-07 Tuesday\", replace: [\"2024-05-08 Wednesday\"] },\n  { find: \"14\u00f77=\", replace: [\"79\u00f76=\"] },\n  { find: \"91\u00f77=\", replace: [\"14\u00f76=\"] },\n  { find: \"37\u00f74=\", replace: [\"57\u00f77=\"] },\n  { find: \"97\u00f75=\", replace: [\"26\u00f73=\"] },\n  { find: \"27\u00f78=\", replace: [\"76\u00f77=\"] },\n  { find: \"28\u00f79=\", replace: [\"60\u00f78=\"] },\n  { find: \"87\u00f78=\", replace: [\"27\u00f74=\"] },\n  { find: \"84\u00f78=\", replace: [\"66\u00f73=\"] },\n  { find: \"60\u00f72=\", replace: [\"20\u00f72=\"] },\n  { find: \"77\u00f72=\", replace: [\"78\u00f73=\"] },\n  { find: \"84\u00f74=\", replace: [\"89\u00f79=\"] },\n  { find: \"53\u00f78=\", replace: [\"29\u00f74=\"] },\n  { find: \"88\u00f76=\", replace: [\"18\u00f72=\"] },\n  { find: \"33\u00f76=\", replace: [\"12\u00f73=\"] },\n  { find: \"13\u00f78=\", replace: [\"12\u00f74=\"] },\n  { find: \"82\u00f79=\", replace: [\"19\u00f73=\"] },\n  { find: \"87\u00f76=\", replace: [\"64\u00f74=\"] },\n  { find: \"69\u00f78=\", replace: [\"97\u00f76=\"] },\n  { find: \"36\u00f75=\", replace: [\"54\u00f73=\", \"33\u00f76=\"] },\n  { find: \"50\u00f79=\", replace: [\"58\u00f77=\"] },\n  { find: \"83\u00f74=\", replace: [\"66\u00f74=\"] },\n  { find: \"20\u00f75=\", replace: [\"10\u00f78=\"] },\n  { find: \"39\u00f72=\", replace: [\"74\u00f74=\"] },\n  { find: \"48\u00f77=\", replace: [\"73\u00f76=\"] },\n];\n\nconst body = context.document.body;\nconst resultSets = [];\nfor (const edit of edits) {\n  const results = body.search(edit.find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  resultSets.push(results);\n}\nawait context.sync();\n\nfor (let i = 0; i < edits.length; i++) {\n  const edit = edits[i];\n  const items = resultSets[i].items;\n  for (let j = 0; j < edit.replace.length; j++) {\n    items[j].insertText(edit.replace[j], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division problems in the practice\n# worksheet table. Each call below searches the whole document content\n# for the exact old text and replaces just that single occurrence\n# (wdReplaceOne), so the one duplicated \"36\u00f75=\" cell is handled\n# correctly: the first call consumes/replaces the first occurrence, and\n# the next search (starting fresh from the top of the document) then\n# finds the still-remaining second occurrence.\n\n$d = $word.ActiveDocument\n\nfunction Replace-FirstMatch($findText, $replaceText) {\n    $r = $d.Content\n    $r.Find.ClearFormatting()\n    $r.Find.Replacement.ClearFormatting()\n    $r.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1) | Out-Null\n}\n\nReplace-FirstMatch \"2024-05-07 Tuesday\" \"2024-05-08 Wednesday\"\n\nReplace-FirstMatch \"14\u00f77=\" \"79\u00f76=\"\nReplace-FirstMatch \"91\u00f77=\" \"14\u00f76=\"\nReplace-FirstMatch \"37\u00f74=\" \"57\u00f77=\"\nReplace-FirstMatch \"97\u00f75=\" \"26\u00f73=\"\nReplace-FirstMatch \"27\u00f78=\" \"76\u00f77=\"\n\nReplace-FirstMatch \"28\u00f79=\" \"60\u00f78=\"\nReplace-FirstMatch \"87\u00f78=\" \"27\u00f74=\"\nReplace-FirstMatch \"84\u00f78=\" \"66\u00f73=\"\nReplace-FirstMatch \"60\u00f72=\" \"20\u00f72=\"\nReplace-FirstMatch \"77\u00f72=\" \"78\u00f73=\"\n\nReplace-FirstMatch \"84\u00f74=\" \"89\u00f79=\"\nReplace-FirstMatch \"53\u00f78=\" \"29\u00f74=\"\nReplace-FirstMatch \"88\u00f76=\" \"18\u00f72=\"\nReplace-FirstMatch \"33\u00f76=\" \"12\u00f73=\"\nReplace-FirstMatch \"13\u00f78=\" \"12\u00f74=\"\n\nReplace-FirstMatch \"82\u00f79=\" \"19\u00f73=\"\nReplace-FirstMatch \"87\u00f76=\" \"64\u00f74=\"\nReplace-FirstMatch \"69\u00f78=\" \"97\u00f76=\"\nReplace-FirstMatch \"36\u00f75=\" \"54\u00f73=\"\nReplace-FirstMatch \"50\u00f79=\" \"58\u00f77=\"\n\nReplace-FirstMatch \"36\u00f75=\" \"33\u00f76=\"\nReplace-FirstMatch \"83\u00f74=\" \"66\u00f74=\"\nReplace-FirstMatch \"20\u00f75=\" \"10\u00f78=\"\nReplace-FirstMatch \"39\u00f72=\" \"74\u00f74=\"\nReplace-FirstMatch \"48\u00f77=\" \"73\u00f76=\"\n"}
